# LOM3099.xlsx update
# - Insert 3 new rows (13-15) for the two additional professors that were
#   added alongside the existing one, turning the single "name" row into a
#   block of three name rows.
# - Rewrite several long-text cells with the updated course description
#   (Objetivos, Programa resumido, Programa, Bibliografia).
# - Adjust column definitions (col A no longer shares a width entry with
#   col B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three blank rows before row 13 (shifts old rows 13-23 down to
#    16-26). This makes room for the two new "Docentes responsaveis" rows.
# ---------------------------------------------------------------------
$ws.Rows("13:15").Insert()

# New row 13 has no value yet -> remove the stray formatted cell Excel
# left behind in column A (the source rows only have B/C content here).
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()

# Copy the B/C formatting from an existing "value" row (row 10) onto the
# three new rows so the shared styles (wrap text, red text in col C, etc.)
# are reused instead of creating new style entries.
$ws.Range("B10").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13:C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the new / updated text content.
# ---------------------------------------------------------------------

# Objetivos: (row 10) now shows the actual objectives paragraph instead of
# the professor's name that had been pasted there by mistake.
$ws.Range("B10").Value = 'Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na estática. Fornecer conhecimentos necessários para cálculo de reações de apoios e de esforços internos em estruturas isostáticas.'
$ws.Range("C10").Value = 'Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na estática. Fornecer conhecimentos necessários para cálculo de reações de apoios e de esforços internos em estruturas isostáticas.'

# Docentes responsaveis: three professors now (rows 13-15).
$ws.Range("B13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("C13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("B14").Value = '5840793 - Sérgio Schneider'
$ws.Range("C14").Value = '5840793 - Sérgio Schneider'
$ws.Range("B15").Value = '7797767 - Viktor Pastoukhov'
$ws.Range("C15").Value = '7797767 - Viktor Pastoukhov'

# Programa resumido: (row 16, previously row 13) gets the actual short
# syllabus text.
$ws.Range("B16").Value = 'Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas.'
$ws.Range("C16").Value = 'Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas.'

# Programa: (row 18, previously row 15) gets the full syllabus text.
$ws.Range("B18").Value = 'Mecânica e suas áreas: Corpos rígidos e corpos deformáveis (sólidos). Terminologia e metodologia básica. Estática de Partículas: Vetores, resultante de várias forças concorrentes, equilíbrio de uma partícula. Estática de Corpos Rígidos: Conceito de corpo rígido. Momento de uma força com relação a um ponto, sistemas equivalentes de forças, momento e binário. Apoios e vínculos. Diagrama de corpo livre. Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações estaticamente indeterminadas e vínculos parciais. Equilíbrio de um corpo rígido em 3D. Análise de Estruturas: análise do equilíbrio de estruturas, ação de múltiplas forças, forças internas, terceira Lei de Newton. Treliças: método dos nós, método das seções. Estruturas e Máquinas: transmissão e modificação de forças. Esforços internos em pórticos, vigas, cabos e eixos de transmissão.'
$ws.Range("C18").Value = 'Mecânica e suas áreas: Corpos rígidos e corpos deformáveis (sólidos). Terminologia e metodologia básica. Estática de Partículas: Vetores, resultante de várias forças concorrentes, equilíbrio de uma partícula. Estática de Corpos Rígidos: Conceito de corpo rígido. Momento de uma força com relação a um ponto, sistemas equivalentes de forças, momento e binário. Apoios e vínculos. Diagrama de corpo livre. Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações estaticamente indeterminadas e vínculos parciais. Equilíbrio de um corpo rígido em 3D. Análise de Estruturas: análise do equilíbrio de estruturas, ação de múltiplas forças, forças internas, terceira Lei de Newton. Treliças: método dos nós, método das seções. Estruturas e Máquinas: transmissão e modificação de forças. Esforços internos em pórticos, vigas, cabos e eixos de transmissão.'

# Método: (row 21, previously row 18) now shows the actual evaluation
# method description instead of the third professor's name.
$ws.Range("B21").Value = 'Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários'
$ws.Range("C21").Value = 'Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários'

# Critério: (row 22, previously row 19) now shows the grading criteria
# text that used to be one row up (under Método:).
$ws.Range("B22").Value = 'Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R).'
$ws.Range("C22").Value = 'Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R).'

# Norma de recuperação: (row 23, previously row 20) now shows the
# recovery-grade norm text that used to be one row up (under Critério:).
$ws.Range("B23").Value = ': A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.'
$ws.Range("C23").Value = ': A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0.'

# Bibliografia: (row 24, previously row 21) gets the actual bibliography
# text (it used to incorrectly hold the "Norma de recuperação" text).
$ws.Range("B24").Value = '1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p.'
$ws.Range("C24").Value = '1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p.'

# ---------------------------------------------------------------------
# 3. Column definitions: column A no longer shares its width entry with
#    column B (column B keeps its own, already-present, 60.71 width).
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 29.83
